# Update the "Förändrad" (Changed) date column (C) for rows 2-24
# from serial date 45214 to 45215 (i.e. bump the date by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45214) {
        $cell.Value2 = 45215
    }
}
